# Reverse the order of the comma-separated "Recorded By" entries in column G
# for every data row on the active sheet (Session Analysis Results).
# Rows whose G value has only a single name (no comma) are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)   # column G = 7
    $val = $cell.Value2

    if ($val -ne $null -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ",\s*"

        $reversed = @()
        for ($i = $parts.Count - 1; $i -ge 0; $i--) {
            $reversed += $parts[$i].Trim()
        }

        $newVal = [string]::Join(", ", $reversed)
        $cell.Value2 = $newVal
    }
}
